# Applies the vic-key-outbreaks cluster-list update:
#  - refreshes the "Cluster name" / "Active cases" rows (A2:B49) to the new outbreak list
#  - several clusters were removed, several added, several renamed, and case counts updated
#  - the table grows from 45 data rows (A1:B46) to 48 data rows (A1:B49)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clusterNames = @(
    '12 Sutton Street Apartment Complex North Melbourne',
    '139 Highett St Apartment Complex Richmond',
    '3153 Sacred Heart Community St Kilda Tier 1A',
    '3528 Ottoman Village Aged Care Broadmeadows',
    '3600 Belvedere Aged Care Noble Park',
    '3612 BlueCross Glengowrie',
    '3652 Regis Aged Care Dandenong North',
    '3824 Estia Health South Morang',
    '3961 Heritage Water Gardens Aged Care FacilitySydenham',
    'ACFS Port Logistics Altona',
    'Alfred Health The Alfred Hospital Melbourne',
    'Armstrong Creek School Armstrong Creek',
    'Aspect Autism Spectrum Australia Disability Service Heatherton',
    'Berwick Fields Primary School Berwick',
    'Berwick Lodge Primary School Berwick',
    'Bubup Womindjeka Family and Children''s Centre Port Melbourne',
    'CREST Children''s Sanctuary Dandenong',
    'Clifton Hill Primary School Clifton Hill',
    'Elements Childcare Warralily Armstrong Creek',
    'Horace Petty Public Housing Estate',
    'Inghams Enterprises Somerville',
    'JBS Australia Brooklyn',
    'KingKids Early Learning Centre and Kindergarten Hallam',
    'Kmart Distribution Centre Truganina',
    'Lilydale Motor Inn Lilydale',
    'Lowanna College Newborough',
    'McQuinns Gym Bendigo',
    'Melbourne Metropolitan Remand Centre Ravenhall',
    'Metcash Limited Distribution Centre Laverton North',
    'Monash Health Casey Hospital Emergency Department Tier 1B',
    'Nido Early School Wyndham Vale',
    'Rosewood Downs Special AccommodationHome Dandenong',
    'Saint Augustines Primary School Wodonga',
    'St Mary''s Primary School Swan Hill',
    'St Thereses Primary School Kennington',
    'St Vincents Hospital Emergency Department Melbourne',
    'St. Brendans Catholic Primary School Lakes Entrance',
    'TUROSI PTY LTD Thomastown',
    'The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1A',
    'Vizzarri Farms Koo Wee Rup',
    'Werribee Mercy Hospital Emergency Department',
    'Werribee Mercy Hospital Werribee',
    'Western Health Sunshine Hospital Emergency Department',
    'Wodonga Cemetery Wodonga',
    'Wodonga Senior Secondary College Wodonga',
    'Wodonga South Primary School Wodonga',
    'Woodend Primary School Woodend',
    'Yooralla Disability Residential Care AlfriedaStreet St Albans'
)

$activeCases = @(
    11,
    11,
    11,
    25,
    13,
    19,
    25,
    58,
    12,
    11,
    14,
    13,
    12,
    13,
    20,
    11,
    11,
    15,
    24,
    10,
    10,
    13,
    11,
    10,
    12,
    32,
    22,
    10,
    18,
    12,
    14,
    12,
    10,
    14,
    13,
    25,
    10,
    14,
    10,
    26,
    34,
    10,
    18,
    39,
    13,
    26,
    10,
    11
)

for ($i = 0; $i -lt $clusterNames.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $clusterNames[$i]
    $ws.Cells.Item($row, 2).Value = $activeCases[$i]
}
